# Generate Report for Handback
#
# 3a249cac-9706-49fa-ac47-b589d2ad5e74 finished handback (status flips from
# "Ready for handoff" to "Handed back: in sync with en-US", matching
# 4c675ee6-6f80-4b42-b109-b0a342d14def) and the rows get re-sorted so row 2
# now describes 3a249cac... and row 3 describes 4c675ee6... on every sheet.
# The hyperlink targets stay anchored to their original cell position /
# relationship id -- only the displayed text (and the underlying cell
# value) moves, so we recreate each sheet's hyperlinks from scratch
# (Delete + Add) using the original target URLs, in the same left-to-right,
# top-to-bottom order, pointed at the new cells.

$wb = $excel.ActiveWorkbook
$missing = [System.Reflection.Missing]::Value

$doneStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("B2").Value2 = $doneStatus
$ws1.Range("C2").Value2 = $doneStatus
$ws1.Range("D2").Value2 = "2016-42-20 16:42:59"

$ws1.Range("B3").Value2 = $doneStatus
$ws1.Range("C3").Value2 = $doneStatus
$ws1.Range("D3").Value2 = "2016-41-20 16:41:38"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/537a4a0230174a94af7876801e1bce930c8eefaf/e2e/4c675ee6-6f80-4b42-b109-b0a342d14def.md", $missing, $missing, "3a249cac-9706-49fa-ac47-b589d2ad5e74.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e5571a9f21a071a4a7ead71399b94b06d9556519/e2e/3a249cac-9706-49fa-ac47-b589d2ad5e74.md", $missing, $missing, "4c675ee6-6f80-4b42-b109-b0a342d14def.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("C2").Value2 = $doneStatus
$ws2.Range("E2").Value2 = "2016-03-20 16:42:56"
$ws2.Range("H2").Value2 = "2016-03-20 16:43:22"
$ws2.Range("I2").Value2 = "Include"

$ws2.Range("C3").Value2 = $doneStatus
$ws2.Range("E3").Value2 = "2016-03-20 16:41:35"
$ws2.Range("H3").Value2 = "2016-03-20 16:42:03"
$ws2.Range("I3").Value2 = "Include"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/537a4a0230174a94af7876801e1bce930c8eefaf/e2e/4c675ee6-6f80-4b42-b109-b0a342d14def.md", $missing, $missing, "3a249cac-9706-49fa-ac47-b589d2ad5e74.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/537a4a0230174a94af7876801e1bce930c8eefaf/e2e/4c675ee6-6f80-4b42-b109-b0a342d14def.md", $missing, $missing, ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a66f496d3b08e8c8a07315245b35fab54b1ff713/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4c675ee6-6f80-4b42-b109-b0a342d14def.2a2d07357f65937f003b078c247658e28442398c.zh-cn.xlf", $missing, $missing, "3a249cac-9706-49fa-ac47-b589d2ad5e74.0a275be988e59c45fffa6eb7b4ae8966e49f1301.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/4495e0d6c634b7f8bfd7908ee76d108833759b54/e2e/4c675ee6-6f80-4b42-b109-b0a342d14def.md", $missing, $missing, "3a249cac-9706-49fa-ac47-b589d2ad5e74.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0430b4a6558d23d75c33a5d664bfb2a48a36c8df/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4c675ee6-6f80-4b42-b109-b0a342d14def.2a2d07357f65937f003b078c247658e28442398c.zh-cn.xlf", $missing, $missing, "3a249cac-9706-49fa-ac47-b589d2ad5e74.0a275be988e59c45fffa6eb7b4ae8966e49f1301.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e5571a9f21a071a4a7ead71399b94b06d9556519/e2e/3a249cac-9706-49fa-ac47-b589d2ad5e74.md", $missing, $missing, "4c675ee6-6f80-4b42-b109-b0a342d14def.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/e5571a9f21a071a4a7ead71399b94b06d9556519/e2e/3a249cac-9706-49fa-ac47-b589d2ad5e74.md", $missing, $missing, ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3ddfdf98099f5a97d7b3a767dda9800f18491eaf/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3a249cac-9706-49fa-ac47-b589d2ad5e74.0a275be988e59c45fffa6eb7b4ae8966e49f1301.zh-cn.xlf", $missing, $missing, "4c675ee6-6f80-4b42-b109-b0a342d14def.2a2d07357f65937f003b078c247658e28442398c.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/4495e0d6c634b7f8bfd7908ee76d108833759b54/e2e/3a249cac-9706-49fa-ac47-b589d2ad5e74.md", $missing, $missing, "4c675ee6-6f80-4b42-b109-b0a342d14def.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0430b4a6558d23d75c33a5d664bfb2a48a36c8df/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3a249cac-9706-49fa-ac47-b589d2ad5e74.0a275be988e59c45fffa6eb7b4ae8966e49f1301.zh-cn.xlf", $missing, $missing, "4c675ee6-6f80-4b42-b109-b0a342d14def.2a2d07357f65937f003b078c247658e28442398c.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("C2").Value2 = $doneStatus
$ws3.Range("I2").Value2 = "Include"

$ws3.Range("C3").Value2 = $doneStatus
$ws3.Range("H3").Value2 = "2016-03-20 16:42:10"
$ws3.Range("I3").Value2 = "Include"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/537a4a0230174a94af7876801e1bce930c8eefaf/e2e/4c675ee6-6f80-4b42-b109-b0a342d14def.md", $missing, $missing, "3a249cac-9706-49fa-ac47-b589d2ad5e74.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/537a4a0230174a94af7876801e1bce930c8eefaf/e2e/4c675ee6-6f80-4b42-b109-b0a342d14def.md", $missing, $missing, ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4cdbd323bf5653161489a9f02c5e43d12f1dd4e3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4c675ee6-6f80-4b42-b109-b0a342d14def.2a2d07357f65937f003b078c247658e28442398c.de-de.xlf", $missing, $missing, "3a249cac-9706-49fa-ac47-b589d2ad5e74.0a275be988e59c45fffa6eb7b4ae8966e49f1301.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/5892741b638dd7aeacd270a41f61715f86b9e2e7/e2e/4c675ee6-6f80-4b42-b109-b0a342d14def.md", $missing, $missing, "3a249cac-9706-49fa-ac47-b589d2ad5e74.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a11c84a4181d3ab9aca98d04cceb90aadee2263c/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4c675ee6-6f80-4b42-b109-b0a342d14def.2a2d07357f65937f003b078c247658e28442398c.de-de.xlf", $missing, $missing, "3a249cac-9706-49fa-ac47-b589d2ad5e74.0a275be988e59c45fffa6eb7b4ae8966e49f1301.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e5571a9f21a071a4a7ead71399b94b06d9556519/e2e/3a249cac-9706-49fa-ac47-b589d2ad5e74.md", $missing, $missing, "4c675ee6-6f80-4b42-b109-b0a342d14def.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/e5571a9f21a071a4a7ead71399b94b06d9556519/e2e/3a249cac-9706-49fa-ac47-b589d2ad5e74.md", $missing, $missing, ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a8842be2353f7775a76519ec86d780518bf664b0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3a249cac-9706-49fa-ac47-b589d2ad5e74.0a275be988e59c45fffa6eb7b4ae8966e49f1301.de-de.xlf", $missing, $missing, "4c675ee6-6f80-4b42-b109-b0a342d14def.2a2d07357f65937f003b078c247658e28442398c.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/5892741b638dd7aeacd270a41f61715f86b9e2e7/e2e/3a249cac-9706-49fa-ac47-b589d2ad5e74.md", $missing, $missing, "4c675ee6-6f80-4b42-b109-b0a342d14def.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a11c84a4181d3ab9aca98d04cceb90aadee2263c/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3a249cac-9706-49fa-ac47-b589d2ad5e74.0a275be988e59c45fffa6eb7b4ae8966e49f1301.de-de.xlf", $missing, $missing, "4c675ee6-6f80-4b42-b109-b0a342d14def.2a2d07357f65937f003b078c247658e28442398c.de-de.xlf") | Out-Null

Write-Host "Done."
